$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Handle the 179x3= collision first: the pre-existing "179×3=" must become
# "144×7=" before we introduce a *new* "179×3=" (from "718×5="), otherwise
# the blanket replace-all would also catch the freshly written text.
Replace-Text "179×3=" "144×7="

Replace-Text "188×7=" "823×4="
Replace-Text "691×5=" "627×8="
Replace-Text "614×6=" "796×5="
Replace-Text "560×5=" "285×3="
Replace-Text "867×6=" "998×7="
Replace-Text "233×2=" "212×3="
Replace-Text "595×4=" "972×5="
Replace-Text "141×5=" "855×8="
Replace-Text "519×6=" "911×5="
Replace-Text "548×5=" "271×8="
Replace-Text "181×4=" "733×7="
Replace-Text "823×6=" "339×8="
Replace-Text "885×9=" "560×8="
Replace-Text "718×5=" "179×3="
Replace-Text "596×3=" "725×6="
Replace-Text "543×9=" "462×4="
Replace-Text "137×3=" "323×4="
Replace-Text "489×9=" "361×9="
Replace-Text "632×7=" "417×5="
Replace-Text "845×3=" "363×5="
Replace-Text "522×8=" "852×8="
Replace-Text "894×7=" "819×6="
Replace-Text "141×4=" "912×8="
Replace-Text "385×9=" "243×2="
